$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I5").Value  = 95.85211039444827
$ws.Range("I6").Value  = 95.96904472825156
$ws.Range("I7").Value  = 95.95666118772137

$ws.Range("G20").Value = 97.83818854344644
$ws.Range("G21").Value = 97.89039325426559
$ws.Range("G22").Value = 97.82887617419976

$ws.Range("H23").Value = 97.36845174419632
$ws.Range("H24").Value = 97.30369602920673
$ws.Range("H25").Value = 97.29430273867894

$ws.Range("I28").Value = 95.83605922117484
$ws.Range("I29").Value = 95.89474369035577

$ws.Range("G38").Value = 97.78139190006415
$ws.Range("G39").Value = 97.77128583716058

$ws.Range("H40").Value = 97.36150567623156
$ws.Range("H41").Value = 97.35078342405272

$ws.Range("I44").Value = 95.94980278223566
$ws.Range("I45").Value = 96.11252830910225

$ws.Range("G54").Value = 97.69000389849066
$ws.Range("G55").Value = 97.89214066011566

$ws.Range("H56").Value = 97.35126222589246
$ws.Range("H57").Value = 97.21915792135152
